$d = $word.ActiveDocument

# Step 1: Fix the misspelled "BäfrFoods" -> replace it with the first part
# of the corrected text ("BärFoo"), keeping the same run formatting.
$d.Content.Find.Execute("BäfrFoods", $true, $false, $false, $false, $false,
                         $true, 1, $false, "BärFoo", 2)

# Step 2: Locate the end of the just-inserted "BärFoo" text - this is the
# seam where the correction happened.
$r = $d.Content
$r.Find.Execute("BärFoo", $true, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)
$r.Collapse(0)
$seamStart = $r.Start

# Step 3: Insert the remainder of the word ("ds") right after the seam,
# producing a new run.
$r.InsertAfter("ds")

# Step 4: Drop a "_GoBack" bookmark exactly at the seam between the two
# runs - this is what Word itself does at the site of the most recent
# edit. Adding it here also removes the previous "_GoBack" bookmark that
# used to sit in the empty paragraph at the end of the document, since
# Word keeps only a single "_GoBack" bookmark at a time.
$seam = $d.Range($seamStart, $seamStart)
$d.Bookmarks.Add("_GoBack", $seam)
